$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before the existing row 334, shifting the old
# rows 334-340 down to become rows 340-346 (their content is untouched).
$ws.Range("A334:R339").EntireRow.Insert()

# Common (fixed) values shared by every row in this data block.
$marketId = 4
$market = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$categoryId = 100112028
$category = "Sandia"
$variety = "Sin especificar"
$kgOUnidades = 1
$classification = "Hortaliza"

# Fill in the 6 newly inserted rows (334-339) with their data.
$rows = 334, 335, 336, 337, 338, 339
$fechas = 44939, 44939, 44939, 44939, 44939, 44939
$calidades = "Cuarta", "Extra", "Tercera", "Segunda", "Super", "Tercera"
$volumenes = 3000, 2000, 2000, 2000, 1500, 2000
$precioMin = 1500, 4000, 3500, 2500, 5000, 2000
$precioMax = 1500, 4000, 3500, 2500, 5000, 2000
$precioProm = 1500, 4000, 3500, 2500, 5000, 2000
$unidades = "$/unidad", "$/unidad", "$/unidad", "$/unidad", "$/unidad", "$/unidad"
$origenes = "Región de O'Higgins", "Región de O'Higgins", "Región de O'Higgins", "Región de O'Higgins", "Región de O'Higgins", "Región de O'Higgins"
$precioKg = 1500, 4000, 3500, 2500, 5000, 2000

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $marketId
    $ws.Cells.Item($r, 2).Value = $market
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fechas[$i]
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $categoryId
    $ws.Cells.Item($r, 7).Value = $category
    $ws.Cells.Item($r, 8).Value = $variety
    $ws.Cells.Item($r, 9).Value = $calidades[$i]
    $ws.Cells.Item($r, 10).Value = $volumenes[$i]
    $ws.Cells.Item($r, 11).Value = $precioMin[$i]
    $ws.Cells.Item($r, 12).Value = $precioMax[$i]
    $ws.Cells.Item($r, 13).Value = $precioProm[$i]
    $ws.Cells.Item($r, 14).Value = $unidades[$i]
    $ws.Cells.Item($r, 15).Value = $origenes[$i]
    $ws.Cells.Item($r, 16).Value = $precioKg[$i]
    $ws.Cells.Item($r, 17).Value = $kgOUnidades
    $ws.Cells.Item($r, 18).Value = $classification
}
